# Inserts one new weekly price record ahead of the existing "Apio" series
# on the active worksheet. All existing records from row 129 down to the
# previous last row (172) are pushed down by one row (to 130-173), and the
# newly available row 129 is populated with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 129:172 down to 130:173, opening up a blank row 129.
$ws.Rows.Item(129).Insert()

# Populate the new row 129 with the new weekly record.
$ws.Cells.Item(129, 1).Value = 7
$ws.Cells.Item(129, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(129, 3).Value = "Ñuble"
$ws.Cells.Item(129, 4).Value = 44588
$ws.Cells.Item(129, 5).Value = 16
$ws.Cells.Item(129, 6).Value = 100112017
$ws.Cells.Item(129, 7).Value = "Apio"
$ws.Cells.Item(129, 8).Value = "Americana (o)"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 60
$ws.Cells.Item(129, 11).Value = 8000
$ws.Cells.Item(129, 12).Value = 8500
$ws.Cells.Item(129, 13).Value = 8250
$ws.Cells.Item(129, 14).Value = "$/docena de matas"
$ws.Cells.Item(129, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(129, 16).Value = 1375
$ws.Cells.Item(129, 17).Value = 6
$ws.Cells.Item(129, 18).Value = "Hortaliza"
